# Add hydrogen combined cycle as a power plant type (#99)
#
# The MCF sheet previously had a single "hydrogen" row (row 24). This
# renames that row to "hydrogen combustion turbine" and adds a brand new
# row ("hydrogen combined cycle") right below it, re-using the same
# capacity factor as natural gas combined cycle (=B4), matching the other
# "derived" rows further up the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCF")

# Row 24: "hydrogen" -> "hydrogen combustion turbine", with the slightly
# different formatting (explicit black font color + vertically centered)
# used for this new crop of rows.
$ws.Range("A24").Value = "hydrogen combustion turbine"
$ws.Range("A24").Font.Italic = $false
$ws.Range("A24").Font.Color = 0
$ws.Range("A24").VerticalAlignment = -4108  # xlCenter

# New row 25: "hydrogen combined cycle", same capacity factor as natural
# gas combined cycle (B4) and the same formatting as row 24.
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B25").Formula = "=B4"

# Restore the original active sheet/selection so only the MCF sheet's
# own selection changes.
$ws.Range("B26").Select() | Out-Null
$wb.Worksheets.Item("About").Activate() | Out-Null
